# Add new columns I ("I0") and J ("IF") to the sheet, matching the
# existing header style used by column H ("IP"), then fill in the
# per-row data values for rows 2-36.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header cell formatting (bold, centered, bordered) from H1
# onto the two new header cells, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I and J, rows 2 through 36.
$iVals = @(6, 7, 8, 5, 7, 5, 6, 6, 6, 6, 5, 8, 5, 5, 5, 6, 6, 8, 8, 5, 7, 5, 9, 7, 7, 8, 9, 3, 9, 10, 7, 1, 6, 4, 3)
$jVals = @(6, 7, 8, 5, 7, 6, 6, 6, 7, 6, 5, 8, 6, 5, 6, 6, 7, 8, 8, 5, 7, 5, 9, 7, 7, 9, 9, 4, 9, 11, 7, 2, 6, 4, 3)

for ($idx = 0; $idx -lt $iVals.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$idx]
    $ws.Cells.Item($row, 10).Value = $jVals[$idx]
}
